# Update the two-digit multiplication problems to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @("96×74=", "98×59="),
    @("88×27=", "52×13="),
    @("46×18=", "51×63="),
    @("91×77=", "43×25="),
    @("59×84=", "82×32="),
    @("95×71=", "30×58="),
    @("60×59=", "39×52="),
    @("76×58=", "95×18="),
    @("12×43=", "31×30="),
    @("19×97=", "14×16="),
    @("21×19=", "69×28="),
    @("95×51=", "94×19="),
    @("84×14=", "64×16="),
    @("27×58=", "40×88="),
    @("76×82=", "23×13="),
    @("92×96=", "97×22="),
    @("18×13=", "40×69="),
    @("83×85=", "57×25="),
    @("34×50=", "70×59="),
    @("16×66=", "33×73="),
    @("75×29=", "51×84="),
    @("71×97=", "88×26="),
    @("14×87=", "56×84="),
    @("50×50=", "47×96="),
    @("42×66=", "89×30=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
